$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.234.76'
$ws.Range('E2').Value = '  +0.26%  '
$ws.Range('D3').Value = '3.112.34'
$ws.Range('E3').Value = '  -0.10%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '579.84'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '173.65'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.33%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  -0.67%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '6.55'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.70%  '
$ws.Range('E10').Value = '  -1.11%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.479'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.61%  '
$ws.Range('E12').Value = '  -0.90%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '36.83'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.50%  '
$ws.Range('E14').Value = '  -1.71%  '
$ws.Range('D15').Value = '3.628.40'
$ws.Range('E15').Value = '  +0.00%  '
$ws.Range('D16').Value = '67.200.78'
$ws.Range('E16').Value = '  +0.15%  '
$ws.Range('E17').Value = '  -1.52%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '16.77'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +2.54%  '
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').Value = '3.114.63'
$ws.Range('E19').Value = '  +0.14%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '490.86'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.81%  '
$ws.Range('B21').Value = 'Polygon'
$ws.Range('C21').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.705'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.32%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.86'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +3.56%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '83.91'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.77%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.15'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.43%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.30'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -3.19%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.58'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +5.08%  '
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.91'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.40%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.36'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -2.66%  '
$ws.Range('E30').Value = '  -0.48%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '28.41'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -2.15%  '
$ws.Range('E32').Value = '  -1.11%  '
$ws.Range('E33').Value = '  -6.76%  '
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.86'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.53%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.972'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.87%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '47.09'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.95%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.04'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -3.72%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.309'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.16%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.124'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.99%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.47'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -2.51%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '391.46'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.93%  '
$ws.Range('D43').Value = '2.802.20'
$ws.Range('E43').Value = '  -1.46%  '
$ws.Range('E44').Value = '  -8.09%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0351'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -2.15%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '135.06'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.66%  '
$ws.Range('E47').Value = '  +0.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '25.10'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.35%  '
$ws.Range('E49').Value = '  -1.06%  '
$ws.Range('E50').Value = '  -1.19%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.71'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -2.78%  '
